$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.183.05"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.825.20"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.92"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6032"
$ws.Range("E6").Value = "  -3.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07113"
$ws.Range("E8").Value = "  -4.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2812"
$ws.Range("E9").Value = "  -2.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.01"
$ws.Range("E10").Value = "  -3.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07641"
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.820.78"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.778"
$ws.Range("E13").Value = "  -4.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6397"
$ws.Range("E14").Value = "  -5.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009791"
$ws.Range("E15").Value = "  -3.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "79.33"
$ws.Range("E16").Value = "  -3.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.044.42"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.955"
$ws.Range("E18").Value = "  -4.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.173.03"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "230.89"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.029"
$ws.Range("E23").Value = "  -5.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.15"
$ws.Range("E25").Value = "  -2.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.012"
$ws.Range("E26").Value = "  -5.41%  "
$ws.Range("E27").Value = "  -5.84%  "
$ws.Range("E28").Value = "  -4.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06708"
$ws.Range("E29").Value = "  +3.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.451"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.455"
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.815"
$ws.Range("E32").Value = "  -6.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.774"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.133"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.706"
$ws.Range("E35").Value = "  -7.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6563"
$ws.Range("E36").Value = "  -5.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.534"
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.233.21"
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.759"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01756"
$ws.Range("E40").Value = "  -5.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.507"
$ws.Range("E41").Value = "  -3.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9248"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.970.36"
$ws.Range("E44").Value = "  -2.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.06"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "62.96"
$ws.Range("E46").Value = "  -4.04%  "
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.622"
$ws.Range("E48").Value = "  -5.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.555"
$ws.Range("E49").Value = "  -7.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05579"
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.451"
$ws.Range("E51").Value = "  -5.92%  "
